$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 28 - new departure entry (Sunday, Jan 15 - FR2474 to London, EI-HGW)
$ws.Range("A28").Value = 27
$ws.Range("B28").Value = "Sunday, Jan 15"
$ws.Range("C28").Value = "10:25 AM"
$ws.Range("D28").Value = "FR2474"
$ws.Range("E28").Value = "London"
$ws.Range("F28").Value = "(STN)"
$ws.Range("G28").Value = "Ryanair "
$ws.Range("H28").Value = "B38M"
$ws.Range("I28").Value = "(EI-HGW)"
$ws.Range("J28").Value = "10:16 AM"
$ws.Range("L28").Value = "0 hours, -9 minutes"

# Row 29 - new departure entry (Sunday, Jan 15 - LO3994 to Warsaw, SP-LIK)
$ws.Range("A29").Value = 28
$ws.Range("B29").Value = "Sunday, Jan 15"
$ws.Range("C29").Value = "3:00 PM"
$ws.Range("D29").Value = "LO3994"
$ws.Range("E29").Value = "Warsaw"
$ws.Range("F29").Value = "(WAW)"
$ws.Range("G29").Value = "LOT "
$ws.Range("H29").Value = "E75S"
$ws.Range("I29").Value = "(SP-LIK)"
$ws.Range("J29").Value = "3:03 PM"
$ws.Range("L29").Value = "0 hours, 3 minutes"

# K and M columns stay blank for both new rows, but the cells still need to
# exist (matching the blank K/M cells present on every other data row).
# Copy the formatting of an existing blank cell so no extra style gets
# introduced and the cells materialize with the default style.
$ws.Range("K27").Copy() | Out-Null
$ws.Range("K28:K29").PasteSpecial(-4122) | Out-Null

$ws.Range("M27").Copy() | Out-Null
$ws.Range("M28:M29").PasteSpecial(-4122) | Out-Null

$excel.CutCopyMode = 0
